$wb = $excel.ActiveWorkbook

# --- CTHPPatientCard sheet: add rows 4-5 (Spanish leukemia cards) ---
$wsPatient = $wb.Worksheets.Item("CTHPPatientCard")
$wsPatient.Range("D4").Value = "CTHP:2"
$wsPatient.Range("C4").Value = "Leucemia de células pilosas"
$wsPatient.Range("B4").Value = "Tratamiento"
$wsPatient.Range("A4").Value = "/espanol/tipos/leucemia"

$wsPatient.Range("A5").Value = "/espanol/tipos/leucemia"
$wsPatient.Range("C5").Value = "Aspectos generales de la prevención del cáncer (PDQ®)"
$wsPatient.Range("B5").Value = "Causas y prevención"
$wsPatient.Range("D5").Value = "CTHP:3"

# --- CTHPHPCard sheet: add row 4 (Spanish "cabeza y cuello" screening card) ---
$wsHP = $wb.Worksheets.Item("CTHPHPCard")
$wsHP.Range("A4").Value = "/espanol/tipos/cabeza-cuello"
$wsHP.Range("D4").Value = "Detección del cáncer de cavidad oral, faringe y laringe"
$wsHP.Range("E4").Value = "CTHP:5"
$wsHP.Range("C4").Value = "Exámenes de detección"
$wsHP.Range("B4").Value = "HP"

# --- Column width adjustments (bestFit widths grew with the new, longer content) ---
$wsHP.Columns.Item(1).ColumnWidth = 26.833333333333332
$wsHP.Columns.Item(4).ColumnWidth = 48.666666666666664

$wsPatient.Columns.Item(1).ColumnWidth = 22.666666666666668
$wsPatient.Columns.Item(2).ColumnWidth = 18.333333333333332

# --- Selection / active sheet bookkeeping ---
# The patient-card sheet keeps its own selection anchored below its new last row,
# but is no longer the active tab.
$wsPatient.Range("A6").Select() | Out-Null

# The HP-card sheet becomes the active tab, with selection below its new last row.
$wsHP.Activate() | Out-Null
$wsHP.Range("A5").Select() | Out-Null
